$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": add a new day column BW (27-aug) with its 24 hourly
# prices, mirroring the style of the existing BV ("26-aug") column.
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell, same style as the other date headers in row 1.
$wsPrix.Range("BW1").Value = "27-aug"
$wsPrix.Range("BV1").Copy()
$wsPrix.Range("BW1").PasteSpecial(-4122)   # xlPasteFormats

$prixValues = @{
    2  = 97.41
    3  = 90.02
    4  = 86.55
    5  = 80.36
    6  = 79.04000000000001
    7  = 85.90000000000001
    8  = 102.84
    9  = 113.12
    10 = 113.64
    11 = 101.63
    12 = 85.20999999999999
    13 = 77.75
    14 = 63.65
    15 = 49.85
    16 = 50.62
    17 = 50.62
    18 = 63.81
    19 = 78.73
    20 = 89.61
    21 = 99
    22 = 109.1
    23 = 114.32
    24 = 113.63
    25 = 103.57
}

foreach ($row in $prixValues.Keys) {
    $wsPrix.Cells.Item($row, 75).Value = $prixValues[$row]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append the new day row (2025-08-25 / 32.625).
# The date column holds plain text (not real dates) in this workbook, so
# force text entry (via a temporary "@" format) to stop auto-date-parsing,
# then drop back to the default "Normal" style to match the other rows,
# which carry no explicit cell style.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A72").NumberFormat = "@"
$wsGaz.Range("A72").Value = "2025-08-25"
$wsGaz.Range("A72").Style = "Normal"
$wsGaz.Range("B72").Value = 32.625

# ---------------------------------------------------------------------------
# Sheet "CO2": append the new day row (2025-08-25 / 71.52).
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A72").NumberFormat = "@"
$wsCo2.Range("A72").Value = "2025-08-25"
$wsCo2.Range("A72").Style = "Normal"
$wsCo2.Range("B72").Value = 71.52
